# Applies: "added Tokenization and Collapse Rare Category" commit
# Adds three new worksheets after "Text Case": "Remove Stopwords",
# "Collapse Rare Categories" and "Tokenization", each following the same
# Action/Time/Content layout used by the other task sheets in this workbook.

$wb = $excel.ActiveWorkbook

# --- Give the pre-existing "Text Case" sheet the selection state it has
#     after the user moved away from it (G21, no longer the active tab). ---
$textCase = $wb.Worksheets.Item("Text Case")
$textCase.Activate()
[void]$textCase.Range("G21").Select()

# --- Create the three new sheets by copying "Text Case" (keeps its
#     Action/Time/Content column styling s="3" header / s="4" body). ---
$textCase.Copy($null, $textCase)
$stopwords = $wb.Worksheets.Item($textCase.Index + 1)
$stopwords.Name = "Remove Stopwords"

$stopwords.Copy($null, $stopwords)
$collapse = $wb.Worksheets.Item($stopwords.Index + 1)
$collapse.Name = "Collapse Rare Categories"

$collapse.Copy($null, $collapse)
$tokenization = $wb.Worksheets.Item($collapse.Index + 1)
$tokenization.Name = "Tokenization"

# ------------------------------------------------------------------
# Remove Stopwords (A1:C6)
# ------------------------------------------------------------------
$stopwords.Range("A1").Value = "Action"
$stopwords.Range("B1").Value = "Time"
$stopwords.Range("C1").Value = "Content"

$stopwords.Range("A2").Value = "Upload CSV"
$stopwords.Range("B2").Value = "5 min"
$stopwords.Range("C2").Value = "df = pd.read_csv('file.csv')"

$stopwords.Range("A3").Value = "Load Stopwords"
$stopwords.Range("B3").Value = "1 min"
$stopwords.Range("C3").Value = "from nltk.corpus import stopwords"

$stopwords.Range("A4").Value = "Remove Stopwords"
$stopwords.Range("B4").Value = "5 min"
$stopwords.Range("C4").Value = "Loop over text column to filter out stopwords"

$stopwords.Range("A5").Value = "Verify Changes"
$stopwords.Range("B5").Value = "1 min"
$stopwords.Range("C5").Value = "Print some text samples to check"

$stopwords.Range("A6").Value = "Overall"
$stopwords.Range("B6").Value = "12 min"

[void]$stopwords.Range("A1:C6").Select()

# ------------------------------------------------------------------
# Collapse Rare Categories (A1:C7 - has an extra row, so insert one
# more body row before the final "Overall" row so it keeps the bold
# s="3" summary style and the inserted row inherits the s="4" body
# style from the row above it).
# ------------------------------------------------------------------
$collapse.Rows.Item(6).Insert()
$collapse.Rows.Item(6).RowHeight = 17

$collapse.Range("A1").Value = "Action"
$collapse.Range("B1").Value = "Time"
$collapse.Range("C1").Value = "Content"

$collapse.Range("A2").Value = "Upload CSV"
$collapse.Range("B2").Value = "5 min"
$collapse.Range("C2").Value = "df = pd.read_csv('file.csv')"

$collapse.Range("A3").Value = "Value Counts"
$collapse.Range("B3").Value = "2 min"
$collapse.Range("C3").Value = "counts = df['category_column'].value_counts()"

$collapse.Range("A4").Value = "Determine Threshold"
$collapse.Range("B4").Value = "2 min"
$collapse.Range("C4").Value = "Decide on a minimum count for categories to be kept"

$collapse.Range("A5").Value = "Collapse Categories"
$collapse.Range("B5").Value = "3 min"
$collapse.Range("C5").Value = "df['category_column'] = df['category_column'].apply(lambda x: 'Other' if counts[x] < threshold else x)"

$collapse.Range("A6").Value = "Verify Changes"
$collapse.Range("B6").Value = "1 min"
$collapse.Range("C6").Value = "df['category_column'].value_counts()"

$collapse.Range("A7").Value = "Overall"
$collapse.Range("B7").Value = "13 min"

[void]$collapse.Range("A1:C7").Select()

# ------------------------------------------------------------------
# Tokenization (A1:C6)
# ------------------------------------------------------------------
$tokenization.Range("A1").Value = "Action"
$tokenization.Range("B1").Value = "Time"
$tokenization.Range("C1").Value = "Content"

$tokenization.Range("A2").Value = "Upload CSV"
$tokenization.Range("B2").Value = "5 min"
$tokenization.Range("C2").Value = "df = pd.read_csv('file.csv')"

$tokenization.Range("A3").Value = "Load Tokenizer"
$tokenization.Range("B3").Value = "1 min"
$tokenization.Range("C3").Value = "from nltk.tokenize import word_tokenize"

$tokenization.Range("A4").Value = "Tokenize Text"
$tokenization.Range("B4").Value = "5 min"
$tokenization.Range("C4").Value = "df['text_column'].apply(word_tokenize)"

$tokenization.Range("A5").Value = "Verify Changes"
$tokenization.Range("B5").Value = "1 min"
$tokenization.Range("C5").Value = "df['text_column'].head()"

$tokenization.Range("A6").Value = "Overall"
$tokenization.Range("B6").Value = "12 min"

[void]$tokenization.Range("A1:C6").Select()

# --- Final state: "Collapse Rare Categories" is the active/visible tab,
#     matching activeTab="13" (0-based) in the saved workbook view. ---
$collapse.Activate()
[void]$collapse.Range("A1:C7").Select()
